# Add new "Logs" entries documenting work on the stairs-to-fontain feature.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$entries = @(
    @{ Row = 61; Date = 45589; Text = "fix stone drop, add stone gravity, fix enemy spawner" },
    @{ Row = 62; Date = 45602; Text = "start work on stairs to fontain, struggle with tweens" },
    @{ Row = 63; Date = 45603; Text = "add viewport to put items, and spagetti logic to handle hovering and putting items in objects" },
    @{ Row = 64; Date = 45604; Text = "fix fontains, items could be putted and picked by hold action button, can move" }
)

$lastRow = 60

foreach ($entry in $entries) {
    $row = $entry.Row

    # Copy the formatting of the last existing row down to the new row so
    # the new cells keep the same date/text styles already used in the log.
    $ws.Range("A" + $lastRow).Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)

    $ws.Range("B" + $lastRow).Copy()
    $ws.Range("B" + $row).PasteSpecial(-4122)

    $ws.Range("A" + $row).Value = $entry.Date
    $ws.Range("B" + $row).Value = $entry.Text
}

$ws.Range("B67").Select()
